# Carlos Anderson Vargas da Silva's RSVP was confirmed: move his row from
# "Recusados" (was row 4) into "Confirmados" (becomes row 2, right after the
# header, since "Confirmados" was previously empty), with Status changed to
# "Confirmado". Then refresh the "Resumo" totals row accordingly.

$wb = $excel.ActiveWorkbook

$confirmados = $wb.Worksheets.Item("Confirmados")
$recusados   = $wb.Worksheets.Item("Recusados")
$resumo      = $wb.Worksheets.Item("Resumo")

# --- Confirmados: was completely empty, now gets a header row + Carlos's row ---
$confirmados.Range("A1").Value = "#"
$confirmados.Range("B1").Value = "Nome"
$confirmados.Range("C1").Value = "Tipo"
$confirmados.Range("D1").Value = "Idade"
$confirmados.Range("E1").Value = "CodigoConvite"
$confirmados.Range("F1").Value = "Status"
$confirmados.Range("G1").Value = "DataConfirmacao"

$confirmados.Range("A2").Value = 1
$confirmados.Range("B2").Value = "Carlos Anderson Vargas da Silva"
$confirmados.Range("C2").Value = "CRIANCA (40 anos)"
$confirmados.Range("D2").Value = 40
$confirmados.Range("E2").Value = 1240
$confirmados.Range("F2").Value = "Confirmado"
$confirmados.Range("G2").Value = "19/04/2025"

# --- Recusados: remove Carlos's row (was row 4) so only rows 1-3 remain ---
$recusados.Rows.Item(4).Delete()

# --- Resumo: refresh the totals row (row 2) ---
$resumo.Range("A2").Value = 1
$resumo.Range("B2").Value = 2
$resumo.Range("C2").Value = 8
$resumo.Range("D2").Value = 1
$resumo.Range("E2").Value = 1
$resumo.Range("F2").Value = 0
$resumo.Range("G2").Value = 0
